$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3.25
$ws.Range("J3").Value = 2.88
$ws.Range("K3").Value = 2.1
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.03
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.4
$ws.Range("T3").Value = 2.75
$ws.Range("U3").Value = 1.8
$ws.Range("V3").Value = 1.91
$ws.Range("W3").Value = 7.5
$ws.Range("X3").Value = 10
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 29
$ws.Range("AC3").Value = 9.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 34
$ws.Range("AK3").Value = 26
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 251
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 12
$ws.Range("AP3").Value = 23
$ws.Range("AR3").Value = 51
$ws.Range("AS3").Value = 151
$ws.Range("AT3").Value = 2.75
$ws.Range("AU3").Value = 8
$ws.Range("AV3").Value = 51
$ws.Range("AX3").Value = 19
$ws.Range("AY3").Value = 26
$ws.Range("AZ3").Value = 51
$ws.Range("BA3").Value = 81
$ws.Range("BB3").Value = 201
$ws.Range("G4").Value = 2.32
$ws.Range("H4").Value = 3.15
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 2.07
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 7.2
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.85
$ws.Range("S4").Value = 1.42
$ws.Range("T4").Value = 2.67
$ws.Range("U4").Value = 1.7
$ws.Range("V4").Value = 2.05
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 11.75
$ws.Range("Y4").Value = 9
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 28
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 6.2
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 9.5
$ws.Range("AH4").Value = 15.5
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 30
$ws.Range("AM4").Value = 400
$ws.Range("AN4").Value = 4.35
$ws.Range("AO4").Value = 13
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 55
$ws.Range("AR4").Value = 90
$ws.Range("AT4").Value = 2.67
$ws.Range("AU4").Value = 6.9
$ws.Range("AV4").Value = 60
$ws.Range("G5").Value = 7.3
$ws.Range("I5").Value = 1.38
$ws.Range("J5").Value = 6.8
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 1.9
$ws.Range("N5").Value = 7.9
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 3.55
$ws.Range("Q5").Value = 1.78
$ws.Range("R5").Value = 1.98
$ws.Range("S5").Value = 1.37
$ws.Range("T5").Value = 2.85
$ws.Range("U5").Value = 2.02
$ws.Range("V5").Value = 1.7
$ws.Range("W5").Value = 19
$ws.Range("X5").Value = 50
$ws.Range("AA5").Value = 90
$ws.Range("AC5").Value = 7.9
$ws.Range("AD5").Value = 8.5
$ws.Range("AE5").Value = 21
$ws.Range("AG5").Value = 6.4
$ws.Range("AH5").Value = 6.2
$ws.Range("AK5").Value = 11.75
$ws.Range("AL5").Value = 29
$ws.Range("AM5").Value = 900
$ws.Range("AN5").Value = 8.5
$ws.Range("AT5").Value = 2.85
$ws.Range("AV5").Value = 90
$ws.Range("AX5").Value = 6.4
$ws.Range("AY5").Value = 17.5
$ws.Range("AZ5").Value = 18.5
$ws.Range("G6").Value = 1.33
$ws.Range("H6").Value = 4.75
$ws.Range("I6").Value = 7.4
$ws.Range("J6").Value = 1.78
$ws.Range("K6").Value = 2.52
$ws.Range("L6").Value = 6.6
$ws.Range("N6").Value = 9.25
$ws.Range("O6").Value = 1.17
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.52
$ws.Range("R6").Value = 2.37
$ws.Range("S6").Value = 1.29
$ws.Range("T6").Value = 3.3
$ws.Range("U6").Value = 1.82
$ws.Range("V6").Value = 1.9
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 7.2
$ws.Range("Y6").Value = 8.5
$ws.Range("Z6").Value = 8.75
$ws.Range("AB6").Value = 23
$ws.Range("AC6").Value = 9.25
$ws.Range("AD6").Value = 10
$ws.Range("AE6").Value = 19
$ws.Range("AF6").Value = 75
$ws.Range("AG6").Value = 24
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 24
$ws.Range("AJ6").Value = 175
$ws.Range("AL6").Value = 65
$ws.Range("AM6").Value = 500
$ws.Range("AO6").Value = 5.9
$ws.Range("AP6").Value = 15
$ws.Range("AQ6").Value = 15.5
$ws.Range("AR6").Value = 40
$ws.Range("AT6").Value = 3.3
$ws.Range("AU6").Value = 8
$ws.Range("AV6").Value = 65
$ws.Range("AW6").Value = 8.75
$ws.Range("AX6").Value = 40
$ws.Range("BA6").Value = 250
$ws.Range("BB6").Value = 450
$ws.Range("G7").Value = 5.4
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 1.57
$ws.Range("J7").Value = 5.5
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 2.18
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 7.2
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 1.91
$ws.Range("R7").Value = 1.83
$ws.Range("S7").Value = 1.42
$ws.Range("T7").Value = 2.67
$ws.Range("U7").Value = 1.93
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 14.5
$ws.Range("X7").Value = 35
$ws.Range("Y7").Value = 17.5
$ws.Range("Z7").Value = 110
$ws.Range("AA7").Value = 60
$ws.Range("AB7").Value = 60
$ws.Range("AC7").Value = 7.2
$ws.Range("AD7").Value = 7.2
$ws.Range("AE7").Value = 17
$ws.Range("AG7").Value = 6.2
$ws.Range("AH7").Value = 7
$ws.Range("AI7").Value = 8
$ws.Range("AJ7").Value = 11.25
$ws.Range("AK7").Value = 13
$ws.Range("AL7").Value = 28
$ws.Range("AN7").Value = 7
$ws.Range("AO7").Value = 32
$ws.Range("AP7").Value = 37
$ws.Range("AQ7").Value = 200
$ws.Range("AR7").Value = 250
$ws.Range("AS7").Value = 500
$ws.Range("AT7").Value = 2.67
$ws.Range("AU7").Value = 7.9
$ws.Range("AV7").Value = 80
$ws.Range("AW7").Value = 3.35
$ws.Range("AX7").Value = 7.9
$ws.Range("AY7").Value = 18.5
$ws.Range("AZ7").Value = 26
$ws.Range("BA7").Value = 65
$ws.Range("BB7").Value = 300
